# chore(results): Auto-update draw results on excel 2025-12-16T17:46:08Z
#
# Appends the latest Pick 4 draw result as a new row at the bottom of the
# "Results" sheet, mirroring the existing rows' layout exactly:
#   Date | Game | Phase | Result | InsertedAt

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Next empty row right after the existing data block (row 90 -> row 91).
$lastRow = $ws.UsedRange.Rows.Count
$newRow  = $lastRow + 1

# Column A..E values for the new draw, in sheet order.
$rowValues = @(
    "2025-12-16",                           # A: Date
    "Pick 4",                               # B: Game
    "251216",                               # C: Phase
    "7-8-6-3",                              # D: Result
    "2025-12-16T21:46:08.196+04:00"         # E: InsertedAt
)

for ($i = 0; $i -lt $rowValues.Length; $i++) {
    $col   = $i + 1
    $cell  = $ws.Cells.Item($newRow, $col)
    $value = $rowValues[$i]

    # Values such as "2025-12-16" or "251216" look like a date/number to
    # Excel's input parser, but the sheet stores every column as plain
    # text (same as all the prior rows, which are all t="str"). Forcing
    # the Text number format right before assigning the value keeps the
    # literal string instead of it getting auto-converted to a date
    # serial / number; the format is cleared again afterwards so the new
    # cells don't carry any extra formatting compared to the rest of the
    # column.
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()

    # Keep Excel's "number stored as text" warning suppressed for the
    # newly appended cells, consistent with the rest of the column.
    try {
        $cell.Errors.Item(3).Ignore = $true
    } catch {
    }
}
